$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.752.45'
$ws.Range("E2").Value = '  +0.31%  '
$ws.Range("D3").Value = '1.603.76'
$ws.Range("E3").Value = '  +0.37%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '''211.83'
$ws.Range("E5").Value = '  +0.21%  '
$ws.Range("E7").Value = '  +0.19%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("E9").Value = '  +0.02%  '
$ws.Range("E10").Value = '  +0.61%  '
$ws.Range("E11").Value = '  +0.81%  '
$ws.Range("D12").Value = '1.828.70'
$ws.Range("E12").Value = '  +0.36%  '
$ws.Range("D13").Value = '1.605.88'
$ws.Range("E13").Value = '  +0.96%  '
$ws.Range("E14").Value = '  +1.01%  '
$ws.Range("E17").Value = '  -0.73%  '
$ws.Range("E18").Value = '  +2.00%  '
$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").Value = '''1.01'
$ws.Range("E19").Value = '  +0.20%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '''209.20'
$ws.Range("E20").Value = '  -0.32%  '
$ws.Range("E21").Value = '  +0.47%  '
$ws.Range("E22").Value = '  -4.51%  '
$ws.Range("E23").Value = '  +0.84%  '
$ws.Range("E24").Value = '  +0.40%  '
$ws.Range("E25").Value = '  +0.26%  '
$ws.Range("D26").Value = '''7.12'
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("E28").Value = '  +0.16%  '
$ws.Range("D30").Value = '''1.15'
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("E31").Value = '  +0.72%  '
$ws.Range("E32").Value = '  +0.45%  '
$ws.Range("D33").Value = '1.286.74'
$ws.Range("E33").Value = '  -0.01%  '
$ws.Range("E34").Value = '  +1.58%  '
$ws.Range("E35").Value = '  +16.40%  '
$ws.Range("E36").Value = '  +0.24%  '
$ws.Range("E37").Value = '  -5.14%  '
$ws.Range("E38").Value = '  -0.80%  '
$ws.Range("E40").Value = '  -0.08%  '
$ws.Range("E41").Value = '  +0.17%  '
$ws.Range("D42").Value = '''0.777'
$ws.Range("E42").Value = '  -0.74%  '
$ws.Range("D43").Value = '''62.72'
$ws.Range("E43").Value = '  -0.87%  '
$ws.Range("D44").Value = '1.740.79'
$ws.Range("E44").Value = '  +0.37%  '
$ws.Range("D45").Value = '''90.37'
$ws.Range("E45").Value = '  -0.88%  '
$ws.Range("E46").Value = '  +0.29%  '
$ws.Range("E47").Value = '  +0.81%  '
$ws.Range("E48").Value = '  +0.55%  '
$ws.Range("D49").Value = '''7.59'
$ws.Range("E49").Value = '  +3.23%  '
$ws.Range("E50").Value = '  +0.10%  '
$ws.Range("E51").Value = '  +1.80%  '
